# Refresh the cryptocurrency Price / Volume(1h) snapshot values (columns D
# and E) for every coin row on the sheet, matching the latest scrape.
#
# Price strings that look like plain numbers (e.g. "0.7118") would
# otherwise be auto-coerced to numeric cells by Excel's type inference,
# dropping trailing zeros / switching to scientific notation and losing
# the original text formatting, so those cells are explicitly forced to
# Text format before the assignment (and reset back to the default style
# afterwards so no stray formatting is introduced).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.386.51"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.876.66"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7118"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3120"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07802"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08451"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "1.865.28"
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7138"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "29.387.05"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008237"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "2.118.57"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.794"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.072"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.427"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.331"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.92%  "
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05293"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7452"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.694"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01872"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").Value = "1.229.03"
$ws.Range("E39").Value = "  +4.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.728"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.481"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8944"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "110.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "2.016.15"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.818"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  +4.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.408"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4330"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.38%  "
